# Generate Report for Archive
#
# 1. Every "Status" cell currently reading "Ready for handoff" becomes
#    "In Translation" (shared across the Overview sheet's per-locale status
#    columns and each locale sheet's own Status column).
# 2. The Status column(s) are narrower now that the text is shorter, so
#    their stored column width shrinks to match (Excel auto-fit behaviour).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRange = $wsOverview.Range("E2:F3")
foreach ($cell in $overviewRange.Cells) {
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value2 = $newStatus
    }
}
# Narrow columns E and F to their post-edit auto-fit width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- Locale sheets: zh-cn / de-de own Status column (col C) ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $statusRange = $ws.Range("C2:C3")
    foreach ($cell in $statusRange.Cells) {
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value2 = $newStatus
        }
    }
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
